$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-05-27T07:09:09+00:00"
$meta.Range("B11").Value = "Extension created as part of OncoFAIR,Extension créée dans le cadre d'OncoFAIR"

# --- Elements sheet ---
$el = $wb.Worksheets.Item("Elements")

# treatmentRank row (row 5)
$el.Range("L5").Value = "Rank of the treatment in which this prescribed protocol is included within the reference protocol in which this prescription is included"
$el.Range("AK5").Value = "numeroCure"

# dayRank row (row 10)
$el.Range("L10").Value = "Ranking of the day in which this prescribed protocol falls within the reference protocol in which this prescription falls"
$el.Range("AK10").Value = "numeroJour"

# dateTimeReference row (row 15)
$el.Range("L15").Value = "Reference date and time of this prescribed protocol"
$el.Range("AK15").Value = "dateHeureReference"
